$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row 1 with new columns P and Q ---
# Copy formatting from O1 (last existing header cell) into P1:Q1
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Swap values in columns I/K and M/O, and add new columns P/Q for rows 2-25 ---
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O

    $ws.Cells.Item($r, 9).Value2  = $kVal  # I <- old K
    $ws.Cells.Item($r, 11).Value2 = $iVal  # K <- old I
    $ws.Cells.Item($r, 13).Value2 = $oVal  # M <- old O
    $ws.Cells.Item($r, 15).Value2 = $mVal  # O <- old M

    $ws.Cells.Item($r, 16).Value2 = 2      # P
    $ws.Cells.Item($r, 17).Value2 = 2      # Q
}

$excel.CutCopyMode = 0
